$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.575530586766542
$ws.Range("C2").Value = 0.0162297128589263
$ws.Range("D2").Value = 0.0599250936329588
$ws.Range("E2").Value = 0.946317103620474
$ws.Range("F2").Value = 0.0112359550561798
$ws.Range("G2").Value = 0.942571785268414
$ws.Range("H2").Value = 0.0137328339575531
$ws.Range("I2").Value = 0.727840199750312
$ws.Range("J2").Value = 0.0424469413233458
$ws.Range("K2").Value = 0.0387016229712859
$ws.Range("L2").Value = 0.0374531835205993
$ws.Range("M2").Value = 0.84019975031211
$ws.Range("N2").Value = 0.00998751560549313
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.0324594257178527
$ws.Range("Q2").Value = 0.928838951310861
$ws.Range("R2").Value = 0.00624219725343321
$ws.Range("S2").Value = 0.00374531835205993
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0.00624219725343321
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.0112359550561798
$ws.Range("X2").Value = 0.00749063670411985
$ws.Range("B3").Value = 0.0649188514357054
$ws.Range("C3").Value = 0.347066167290886
$ws.Range("D3").Value = 0.84019975031211
$ws.Range("E3").Value = 0.0436953807740325
$ws.Range("F3").Value = 0.00124843945068664
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.0699126092384519
$ws.Range("I3").Value = 0.0374531835205993
$ws.Range("J3").Value = 0.131086142322097
$ws.Range("K3").Value = 0.920099875156055
$ws.Range("L3").Value = 0.958801498127341
$ws.Range("M3").Value = 0.108614232209738
$ws.Range("N3").Value = 0.920099875156055
$ws.Range("O3").Value = 0.0149812734082397
$ws.Range("P3").Value = 0.00124843945068664
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.852684144818976
$ws.Range("S3").Value = 0.983770287141074
$ws.Range("T3").Value = 0.0561797752808989
$ws.Range("U3").Value = 0.00749063670411985
$ws.Range("V3").Value = 0.0212234706616729
$ws.Range("W3").Value = 0.0324594257178527
$ws.Range("X3").Value = 0.00124843945068664
$ws.Range("B4").Value = 0.35330836454432
$ws.Range("C4").Value = 0.0249687890137328
$ws.Range("D4").Value = 0.00374531835205993
$ws.Range("E4").Value = 0.00374531835205993
$ws.Range("F4").Value = 0.946317103620474
$ws.Range("G4").Value = 0.0536828963795256
$ws.Range("H4").Value = 0.00249687890137328
$ws.Range("I4").Value = 0.0199750312109863
$ws.Range("J4").Value = 0.0536828963795256
$ws.Range("K4").Value = 0.0387016229712859
$ws.Range("L4").Value = 0.00124843945068664
$ws.Range("M4").Value = 0.00624219725343321
$ws.Range("N4").Value = 0.00374531835205993
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.963795255930087
$ws.Range("Q4").Value = 0.00374531835205993
$ws.Range("R4").Value = 0.133583021223471
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0.00249687890137328
$ws.Range("V4").Value = 0.0536828963795256
$ws.Range("W4").Value = 0.953807740324594
$ws.Range("X4").Value = 0.955056179775281
$ws.Range("B5").Value = 0.00624219725343321
$ws.Range("C5").Value = 0.611735330836454
$ws.Range("D5").Value = 0.0936329588014981
$ws.Range("E5").Value = 0.00499375780274657
$ws.Range("F5").Value = 0.0411985018726592
$ws.Range("G5").Value = 0.00374531835205993
$ws.Range("H5").Value = 0.913857677902622
$ws.Range("I5").Value = 0.214731585518102
$ws.Range("J5").Value = 0.772784019975031
$ws.Range("K5").Value = 0.00124843945068664
$ws.Range("L5").Value = 0.00249687890137328
$ws.Range("M5").Value = 0.0449438202247191
$ws.Range("N5").Value = 0.066167290886392
$ws.Range("O5").Value = 0.98501872659176
$ws.Range("P5").Value = 0.00249687890137328
$ws.Range("Q5").Value = 0.0674157303370786
$ws.Range("R5").Value = 0.00749063670411985
$ws.Range("S5").Value = 0.0124843945068664
$ws.Range("T5").Value = 0.943820224719101
$ws.Range("U5").Value = 0.983770287141074
$ws.Range("V5").Value = 0.925093632958802
$ws.Range("W5").Value = 0.00249687890137328
$ws.Range("X5").Value = 0.0362047440699126
Write-Output "Updated B2:X5 with new frequency values."
